$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (F/G columns, and a few B/D corrections) ---

$ws.Range("F465").Value() = 61661
$ws.Range("G465").Value() = 57

$ws.Range("F472").Value() = 52076

$ws.Range("F475").Value() = 36800
$ws.Range("G475").Value() = 30

$ws.Range("F476").Value() = 37459

$ws.Range("F477").Value() = 37334

$ws.Range("F478").Value() = 55267

$ws.Range("F479").Value() = 42758

$ws.Range("G562").Value() = 284

$ws.Range("F567").Value() = 23492

$ws.Range("F583").Value() = 29155
$ws.Range("G583").Value() = 488

$ws.Range("F584").Value() = 13163
$ws.Range("G584").Value() = 248

$ws.Range("F585").Value() = 14845
$ws.Range("G585").Value() = 358

$ws.Range("F586").Value() = 33448
$ws.Range("G586").Value() = 702

$ws.Range("F587").Value() = 28116
$ws.Range("G587").Value() = 552

$ws.Range("F588").Value() = 25212
$ws.Range("G588").Value() = 512

$ws.Range("B589").Value() = 433708
$ws.Range("D589").Value() = 1951
$ws.Range("F589").Value() = 25275
$ws.Range("G589").Value() = 467

$ws.Range("B590").Value() = 435649
$ws.Range("F590").Value() = 28900
$ws.Range("G590").Value() = 579

$ws.Range("B591").Value() = 437384
$ws.Range("F591").Value() = 14341
$ws.Range("G591").Value() = 412

$ws.Range("B592").Value() = 437937
$ws.Range("F592").Value() = 18126
$ws.Range("G592").Value() = 650

# --- Append new rows 593-595 ---

$ws.Range("A593").Value() = 44487
$ws.Range("B593").Value() = 439734
$ws.Range("C593").Value() = 12160
$ws.Range("D593").Value() = 1797
$ws.Range("E593").Value() = 12864
$ws.Range("F593").Value() = 36578
$ws.Range("G593").Value() = 1177

$ws.Range("A594").Value() = 44488
$ws.Range("B594").Value() = 443214
$ws.Range("C594").Value() = 16199
$ws.Range("D594").Value() = 3480
$ws.Range("E594").Value() = 12872
$ws.Range("F594").Value() = 28986
$ws.Range("G594").Value() = 807

$ws.Range("A595").Value() = 44489
$ws.Range("B595").Value() = 446305
$ws.Range("C595").Value() = 14273
$ws.Range("D595").Value() = 3091
$ws.Range("E595").Value() = 12883
$ws.Range("F595").Value() = 19618
$ws.Range("G595").Value() = 612
